# Phase-two title wording update on slide 1:
#   "Spatio-temporal Remote Sensing, automated in-situ IoT sensors & ANN to
#    monitor and predict HABs and cyanotoxins."
# becomes
#   "Spatio-temporal Remote Sensing & automated in-situ IoT sensors to
#    monitor and predict HABs and cyanotoxins."
#
# The title run that used to read "sensors & ANN to monitor ..." loses the
# "& ANN" and the run that used to read "Spatio-temporal Remote Sensing,
# automated " gets its comma swapped for an ampersand - both edits split
# their host run into separate <a:r> runs (PowerPoint keeps the untouched
# head/tail text as its own run when you retype a chunk in the middle).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# --- Edit the "sensors & ANN to " run first (it sits further right in the
#     paragraph), so the character offsets used for the first run below are
#     not disturbed by a change in text length.
$sensors = $tr.Characters(55, 17)
if ($sensors.Text -ne "sensors & ANN to ") {
    throw "unexpected text at sensors offset: [$($sensors.Text)]"
}
$sensors.Text = "sensors to "

# --- Edit "Remote Sensing, " -> "Remote Sensing & " inside the title run.
$remote = $tr.Characters(17, 16)
if ($remote.Text -ne "Remote Sensing, ") {
    throw "unexpected text at remote-sensing offset: [$($remote.Text)]"
}
$remote.Text = "Remote Sensing & "
